$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Existing rows 2-7: update values in place
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 149.656361
$ws.Range("H2").Value = 448.969083
$ws.Range("I2").Value = 0.5921360794347563
$ws.Range("J2").Value = 0.5921360794347564
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 22.906497
$ws.Range("N2").Value = 68.719491
$ws.Range("O2").Value = 0.9446038650914245
$ws.Range("P2").Value = 0.9446038650914245
$ws.Range("Q2").Value = 3428.102984277417
$ws.Range("R2").Value = 30852.92685849675
$ws.Range("S2").Value = 0.5593340292941535
$ws.Range("T2").Value = 0.5593340292941537

$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 149.656361
$ws.Range("H3").Value = 448.969083
$ws.Range("I3").Value = 0.5921360794347563
$ws.Range("J3").Value = 0.5921360794347564
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.1329193333333333
$ws.Range("N3").Value = 0.3987579999999999
$ws.Range("O3").Value = 0.005481244732096839
$ws.Range("P3").Value = 0.005481244732096839
$ws.Range("Q3").Value = 19.89222373321266
$ws.Range("R3").Value = 179.030013598914
$ws.Range("S3").Value = 0.003245642766086234
$ws.Range("T3").Value = 0.003245642766086234

$ws.Range("A4").Value = "ECs"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 149.656361
$ws.Range("H4").Value = 448.969083
$ws.Range("I4").Value = 0.5921360794347563
$ws.Range("J4").Value = 0.5921360794347564
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.210428333333333
$ws.Range("N4").Value = 3.631285
$ws.Range("O4").Value = 0.04991489017647865
$ws.Range("P4").Value = 0.04991489017647865
$ws.Range("Q4").Value = 181.1482996179616
$ws.Range("R4").Value = 1630.334696561655
$ws.Range("S4").Value = 0.0295564073745165
$ws.Range("T4").Value = 0.0295564073745165

$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 52.73412466666667
$ws.Range("H5").Value = 158.202374
$ws.Range("I5").Value = 0.208649853730866
$ws.Range("J5").Value = 0.208649853730866
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 22.906497
$ws.Range("N5").Value = 68.719491
$ws.Range("O5").Value = 0.9446038650914245
$ws.Range("P5").Value = 0.9446038650914245
$ws.Range("Q5").Value = 1207.954068474626
$ws.Range("R5").Value = 10871.58661627163
$ws.Range("S5").Value = 0.1970914582849364
$ws.Range("T5").Value = 0.1970914582849364

$ws.Range("A6").Value = "FAPs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 52.73412466666667
$ws.Range("H6").Value = 158.202374
$ws.Range("I6").Value = 0.208649853730866
$ws.Range("J6").Value = 0.208649853730866
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1329193333333333
$ws.Range("N6").Value = 0.3987579999999999
$ws.Range("O6").Value = 0.005481244732096839
$ws.Range("P6").Value = 0.005481244732096839
$ws.Range("Q6").Value = 7.009384694610221
$ws.Range("R6").Value = 63.08446225149198
$ws.Range("S6").Value = 0.001143660911615085
$ws.Range("T6").Value = 0.001143660911615085

$ws.Range("A7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 52.73412466666667
$ws.Range("H7").Value = 158.202374
$ws.Range("I7").Value = 0.208649853730866
$ws.Range("J7").Value = 0.208649853730866
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.210428333333333
$ws.Range("N7").Value = 3.631285
$ws.Range("O7").Value = 0.04991489017647865
$ws.Range("P7").Value = 0.04991489017647865
$ws.Range("Q7").Value = 63.83087863006556
$ws.Range("R7").Value = 574.47790767059
$ws.Range("S7").Value = 0.01041473453431451
$ws.Range("T7").Value = 0.01041473453431451

# New rows 8-10 (sCs sender row x 3 targets)
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Thbs1"
$ws.Range("C8").Value = "Itga4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 50.34932566666667
$ws.Range("H8").Value = 151.047977
$ws.Range("I8").Value = 0.1992140668343777
$ws.Range("J8").Value = 0.1992140668343777
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 22.906497
$ws.Range("N8").Value = 68.719491
$ws.Range("O8").Value = 0.9446038650914245
$ws.Range("P8").Value = 0.9446038650914245
$ws.Range("Q8").Value = 1153.326677335523
$ws.Range("R8").Value = 10379.94009601971
$ws.Range("S8").Value = 0.1881783775123345
$ws.Range("T8").Value = 0.1881783775123345

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Thbs1"
$ws.Range("C9").Value = "Itga4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 50.34932566666667
$ws.Range("H9").Value = 151.047977
$ws.Range("I9").Value = 0.1992140668343777
$ws.Range("J9").Value = 0.1992140668343777
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1329193333333333
$ws.Range("N9").Value = 0.3987579999999999
$ws.Range("O9").Value = 0.005481244732096839
$ws.Range("P9").Value = 0.005481244732096839
$ws.Range("Q9").Value = 6.69239880139622
$ws.Range("R9").Value = 60.23158921256599
$ws.Range("S9").Value = 0.00109194105439552
$ws.Range("T9").Value = 0.001091941054395521

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Thbs1"
$ws.Range("C10").Value = "Itga4"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 50.34932566666667
$ws.Range("H10").Value = 151.047977
$ws.Range("I10").Value = 0.1992140668343777
$ws.Range("J10").Value = 0.1992140668343777
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.210428333333333
$ws.Range("N10").Value = 3.631285
$ws.Range("O10").Value = 0.04991489017647865
$ws.Range("P10").Value = 0.04991489017647865
$ws.Range("Q10").Value = 60.94425035116056
$ws.Range("R10").Value = 548.4982531604451
$ws.Range("S10").Value = 0.009943748267647641
$ws.Range("T10").Value = 0.009943748267647641
